$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "Senior Five" ----
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("D2").Value = 20.0
$ws1.Range("F2").Value = 44.0
$ws1.Range("G2").Value = 35.0
$ws1.Range("H2").Value = 57.0
$ws1.Range("I2").Value = 36.0

$ws1.Range("F4").Value = 58.0
$ws1.Range("G4").Value = 32.0
$ws1.Range("H4").Value = 76.0
$ws1.Range("I4").Value = 49.0

$ws1.Range("D6").Value = 20.0
$ws1.Range("F6").Value = 26.0
$ws1.Range("G6").Value = 41.0
$ws1.Range("H6").Value = 40.0
$ws1.Range("I6").Value = 20.0

$ws1.Range("D9").Value = 50.0
$ws1.Range("E9").Value = 46.0
$ws1.Range("F9").Value = 68.0
$ws1.Range("G9").Value = 49.0
$ws1.Range("H9").Value = 77.0
$ws1.Range("I9").Value = 69.0

$ws1.Range("D11").Value = 50.0
$ws1.Range("F11").Value = 62.0
$ws1.Range("G11").Value = 55.0

$ws1.Range("D12").Value = 24.0
$ws1.Range("G12").Value = 20.0

$ws1.Range("D13").Value = 56.0
$ws1.Range("E13").Value = 34.0
$ws1.Range("F13").Value = 42.0
$ws1.Range("H13").Value = 79.0
$ws1.Range("I13").Value = 60.0

$ws1.Range("D16").Value = 26.0
$ws1.Range("E16").Value = 38.0
$ws1.Range("F16").Value = 24.0
$ws1.Range("G16").Value = 45.0
$ws1.Range("H16").Value = 44.0
$ws1.Range("I16").Value = 33.0

$ws1.Range("F17").Value = 42.0
$ws1.Range("G17").Value = 51.0
$ws1.Range("H17").Value = 71.0
$ws1.Range("I17").Value = 48.0

# New students added at the bottom of the Senior Five roster
$ws1.Range("B18").Value = "ACEN BRENDA FASSY"
$ws1.Range("D18").Value = 26.0
$ws1.Range("F18").Value = 20.0
$ws1.Range("G18").Value = 28.0
$ws1.Range("H18").Value = 41.0
$ws1.Range("I18").Value = 27.0

$ws1.Range("B19").Value = "PILOYA MERCY"
$ws1.Range("G19").Value = 43.0
$ws1.Range("H19").Value = 72.0
$ws1.Range("I19").Value = 43.0

# ---- Sheet 2: "Senior Six" ----
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("E3").Value = 21.0
$ws2.Range("F3").Value = 74.0

$ws2.Range("F4").Value = 40.0

$ws2.Range("F6").Value = 12.0

$ws2.Range("F11").Value = 66.0

$ws2.Range("E14").Value = 55.0
$ws2.Range("F14").Value = 68.0

$ws2.Range("E16").Value = 48.0
$ws2.Range("F16").Value = 30.0

$ws2.Range("E18").Value = 56.0
$ws2.Range("F18").Value = 70.0

$ws2.Range("E19").Value = 44.0
